$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns for the two new importable fields
$ws.Range("H1").Value = "reference_period"
$ws.Range("I1").Value = "remarks"

# Give the original header/data range (A1:G1 and A2:F2) an explicit
# "Normal" cell style - this is what produces the second cellXf record
# (applyFont etc.) that the new file carries for the pre-existing cells,
# while the two freshly added header cells (H1/I1) keep using the
# original default style.
$ws.Range("A1:G1").Style = "Normal"
$ws.Range("A2:F2").Style = "Normal"

# Give the new "reference_period" column a sensible custom width, as in
# the target workbook.
$ws.Columns.Item(8).ColumnWidth = 14.6

# Restore the last active cell/selection recorded in the sheet.
$ws.Range("J11").Select() | Out-Null
